$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-typed cells (Price/Volume columns) keep their original
# text formatting (e.g. trailing zeros) instead of being auto-converted
# to numbers by COM Value assignment.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '69.618.11'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '3.673.27'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '630.37'
$ws.Range('E5').Value = '  -6.31%  '
$ws.Range('D6').Value = '160.38'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').Value = '7.15'
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').Value = '0.441'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('D13').Value = '4.291.04'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Value = '32.54'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').Value = '3.669.44'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '69.681.48'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '15.91'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').Value = '10.33'
$ws.Range('E20').Value = '  +5.32%  '
$ws.Range('D21').Value = '471.46'
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').Value = '0.652'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '79.78'
$ws.Range('D24').Value = '3.815.69'
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('D27').Value = '11.09'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  -4.26%  '
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('E30').Value = '  -4.05%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '1.99'
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '26.67'
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '0.986'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('D35').Value = '6.41'
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('D36').Value = '3.674.38'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').Value = '8.33'
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('D39').Value = '178.83'
$ws.Range('E39').Value = '  +3.67%  '
$ws.Range('E40').Value = '  -4.73%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('D43').Value = '0.0894'
$ws.Range('E43').Value = '  -1.32%  '
$ws.Range('D44').Value = '0.927'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '46.69'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '29.10'
$ws.Range('E46').Value = '  +3.69%  '
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('D48').Value = '7.87'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('E49').Value = '  -4.76%  '
$ws.Range('E50').Value = '  -5.50%  '
$ws.Range('E51').Value = '  -4.94%  '
